$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

$ws.Range("A1").Value = "Datos actualizados a 2 de Septiembre de 2020 a las 20:46"
$ws.Range("B4").Value = 6274487
$ws.Range("C4").Value = 16916
$ws.Range("D4").Value = 3512533
$ws.Range("E4").Value = 2572619
$ws.Range("G4").Value = 435
$ws.Range("H4").Value = 189335
$ws.Range("B6").Value = 3848968
$ws.Range("C6").Value = 82860
$ws.Range("D6").Value = 2967396
$ws.Range("E6").Value = 814086
$ws.Range("G6").Value = 1026
$ws.Range("H6").Value = 67486
$ws.Range("D20").Value = 86963
$ws.Range("E20").Value = 175375
$ws.Range("G20").Value = 25
$ws.Range("H20").Value = 30686
$ws.Range("B21").Value = 273301
$ws.Range("C21").Value = 1596
$ws.Range("D21").Value = 246876
$ws.Range("E21").Value = 19963
$ws.Range("G21").Value = 45
$ws.Range("H21").Value = 6462
$ws.Range("B23").Value = 246893
$ws.Range("C23").Value = 892
$ws.Range("E23").Value = 15708
$ws.Range("B49").Value = 65453
$ws.Range("C49").Value = 1672
$ws.Range("D49").Value = 50357
$ws.Range("E49").Value = 13880
$ws.Range("G49").Value = 32
$ws.Range("H49").Value = 1216
$ws.Range("A53").Value = "Etiopia"
$ws.Range("B53").Value = 54409
$ws.Range("C53").Value = 1105
$ws.Range("D53").Value = 19903
$ws.Range("E53").Value = 33660
$ws.Range("G53").Value = 18
$ws.Range("H53").Value = 846
$ws.Range("A54").Value = "Nigeria"
$ws.Range("B54").Value = 54247
$ws.Range("D54").Value = 42010
$ws.Range("E54").Value = 11214
$ws.Range("H54").Value = 1023
$ws.Range("B70").Value = 29114
$ws.Range("C70").Value = 89
$ws.Range("E70").Value = 3973
$ws.Range("B115").Value = 4668
$ws.Range("C115").Value = 50
$ws.Range("D115").Value = 3585
$ws.Range("E115").Value = 989
$ws.Range("B132").Value = 2898
$ws.Range("C132").Value = 68
$ws.Range("D132").Value = 661
$ws.Range("E132").Value = 2117
$ws.Range("G132").Value = 4
$ws.Range("H132").Value = 120
$ws.Range("A140").Value = "Jordania"
$ws.Range("B140").Value = 2161
$ws.Range("C140").Value = 64
$ws.Range("D140").Value = 1610
$ws.Range("E140").Value = 536
$ws.Range("H140").Value = 15
$ws.Range("A141").Value = "Benin"
$ws.Range("B141").Value = 2145
$ws.Range("C141").Value = 0
$ws.Range("D141").Value = 1738
$ws.Range("E141").Value = 367
$ws.Range("H141").Value = 40
$ws.Range("A142").Value = "Islandia"
$ws.Range("B142").Value = 2121
$ws.Range("C142").Value = 5
$ws.Range("D142").Value = 2016
$ws.Range("E142").Value = 95
$ws.Range("A143").Value = "Aruba"
$ws.Range("B143").Value = 2104
$ws.Range("D143").Value = 857
$ws.Range("E143").Value = 1237
$ws.Range("H143").Value = 10
$ws.Range("B144").Value = 2029
$ws.Range("C144").Value = 1
$ws.Range("D144").Value = 1600
$ws.Range("E144").Value = 358
$ws.Range("D158").Value = 1162
$ws.Range("E158").Value = 61
$ws.Range("B160").Value = 1199
$ws.Range("C160").Value = 15
$ws.Range("D160").Value = 909
$ws.Range("E160").Value = 237
